$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "Giovanni Lopes"
$ws.Range("C2").Value = 19
$ws.Range("D2").Value = "Masculino"
$ws.Range("E2").Value = "Pardo"
$ws.Range("F2").Value = "Mariana Sandra Rocha Ribeiro Lopes"
$ws.Range("G2").Value = "Em idade escolar."
$ws.Range("H2").Value = "Empregado: Trabalhador doméstico (sem CLT)"
$ws.Range("I2").Value = "311.492.220-18"
$ws.Range("J2").Value = "903704926"

# Row 3
$ws.Range("B3").Value = "Karla Regina Cavalcanti Rocha"
$ws.Range("C3").Value = 20
$ws.Range("D3").Value = "Feminino"
$ws.Range("E3").Value = "Parda"
$ws.Range("F3").Value = "Jessica Gouveia Cavalcanti Rocha Regina"
$ws.Range("H3").Value = "Empregado: Setor privado (CLT)"
$ws.Range("I3").Value = "709.951.750-39"
$ws.Range("J3").Value = "986764036"

# Row 4
$ws.Range("B4").Value = "Tatiane Barros"
$ws.Range("C4").Value = 23
$ws.Range("D4").Value = "Feminino"
$ws.Range("E4").Value = "Parda"
$ws.Range("F4").Value = "Ana Oliveira Barros"
$ws.Range("H4").Value = "Autonomo: Sem CNPJ"
$ws.Range("I4").Value = "234.431.350-80"
$ws.Range("J4").Value = "139580629"

# Row 5
$ws.Range("B5").Value = "Pedro Gouveia Ribeira"
$ws.Range("C5").Value = 1
$ws.Range("F5").Value = "Rafaela Carolina Gouveia Ribeira"
$ws.Range("I5").Value = "131.169.770-54"
$ws.Range("J5").Value = "503038146"

# Row 6
$ws.Range("B6").Value = "Natália Moraes Pires"
$ws.Range("C6").Value = 19
$ws.Range("D6").Value = "Feminino"
$ws.Range("E6").Value = "Parda"
$ws.Range("F6").Value = "Leticia Costa Pires Moraes"
$ws.Range("I6").Value = "568.828.350-55"
$ws.Range("J6").Value = "445039353"

# Row 7
$ws.Range("B7").Value = "Helena Campos"
$ws.Range("C7").Value = 17
$ws.Range("D7").Value = "Feminino"
$ws.Range("E7").Value = "Branca"
$ws.Range("F7").Value = "Brenda Silva Campos"
$ws.Range("H7").Value = "Empregado: Setor privado (sem CLT)"
$ws.Range("I7").Value = "501.891.060-77"
$ws.Range("J7").Value = "730157861"

# Row 8
$ws.Range("B8").Value = "Laura Nascimento Santana Alves"
$ws.Range("C8").Value = 21
$ws.Range("F8").Value = "Rosana Alves Nascimento"
$ws.Range("G8").Value = "Em idade escolar."
$ws.Range("H8").Value = "Empregado: Setor privado (CLT)"
$ws.Range("I8").Value = "135.803.630-69"
$ws.Range("J8").Value = "266613963"

# Row 9
$ws.Range("B9").Value = "Elaine Alves Machado"
$ws.Range("C9").Value = 38
$ws.Range("F9").Value = "Nathalia Alves Machado"
$ws.Range("G9").Value = "Superior incompleto"
$ws.Range("H9").Value = "Empregado: Setor publico (sem CLT)"
$ws.Range("I9").Value = "149.076.830-06"
$ws.Range("J9").Value = "366229579"

# Row 10
$ws.Range("B10").Value = "Ricardo Gouveia Andrade"
$ws.Range("C10").Value = 9
$ws.Range("D10").Value = "Masculino"
$ws.Range("E10").Value = "Pardo"
$ws.Range("F10").Value = "Luiza Gouveia"
$ws.Range("G10").Value = "Em idade escolar."
$ws.Range("H10").Value = "Fora da força de trabalho"
$ws.Range("I10").Value = "987.971.870-46"
$ws.Range("J10").Value = "455452585"

# Row 11
$ws.Range("B11").Value = "Larissa Cecília Cardoso Martins"
$ws.Range("C11").Value = 3
$ws.Range("D11").Value = "Feminino"
$ws.Range("E11").Value = "Branca"
$ws.Range("F11").Value = "Carolina Cardoso"
$ws.Range("H11").Value = "Fora da força de trabalho"
$ws.Range("I11").Value = "407.675.280-50"
$ws.Range("J11").Value = "674898996"
